$wb = $excel.ActiveWorkbook

# --- Sheet1: update the summary counts for the SLA/release may row ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("H2").Value = 5
$ws1.Range("J2").Value = 5

# --- Text_Summary_REPORT: update the matching narrative text ---
$ws2 = $wb.Worksheets.Item("Text_Summary_REPORT")
$ws2.Range("A3").Value = "5 Test cases designed (fixVersion)."
$ws2.Range("A4").Value = "Story Testing - 2 out of 5 test cases executed."
